# Auto-generated Excel COM-interop script to apply the F/G column updates
# described in the commit diff (gh-pages data refresh at 456a3b4).
$wb = $excel.ActiveWorkbook

# --- sheet1 (Worksheets.Item(1)) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 7846
$ws.Range("F3").Value = 7989
$ws.Range("F5").Value = 43
$ws.Range("F6").Value = 6971
$ws.Range("F7").Value = 3454
$ws.Range("F9").Value = 3767
$ws.Range("F12").Value = 59
$ws.Range("F14").Value = 105
$ws.Range("F15").Value = 490
$ws.Range("F16").Value = 14
$ws.Range("F17").Value = 89
$ws.Range("F18").Value = 340
$ws.Range("F21").Value = 344
$ws.Range("F22").Value = 3982
$ws.Range("F24").Value = 391
$ws.Range("F25").Value = 1017
$ws.Range("F26").Value = 512
$ws.Range("F27").Value = 1574
$ws.Range("F28").Value = 91
$ws.Range("F30").Value = 2886
$ws.Range("F31").Value = 2063
$ws.Range("F32").Value = 46
$ws.Range("F33").Value = 66
$ws.Range("F34").Value = 74
$ws.Range("F35").Value = 106
$ws.Range("F36").Value = 3984
$ws.Range("F37").Value = 391
$ws.Range("F39").Value = 46
$ws.Range("F41").Value = 687
$ws.Range("F42").Value = 102
$ws.Range("F44").Value = 1528
$ws.Range("F47").Value = 585

# --- sheet2 (Worksheets.Item(2)) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("G3").Value = "已停售"
$ws.Range("F4").Value = 428
$ws.Range("F7").Value = 58
$ws.Range("F10").Value = 42
$ws.Range("F16").Value = 480

# --- sheet3 (Worksheets.Item(3)) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 152

# --- sheet4 (Worksheets.Item(4)) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 152
$ws.Range("G4").Value = "已停售"
$ws.Range("F5").Value = 7846
$ws.Range("F6").Value = 7989
$ws.Range("F8").Value = 43
$ws.Range("F9").Value = 6971
$ws.Range("F10").Value = 3454
$ws.Range("F12").Value = 3767
$ws.Range("F15").Value = 59
$ws.Range("F16").Value = 105
$ws.Range("F17").Value = 14
$ws.Range("F19").Value = 89
$ws.Range("F21").Value = 58
$ws.Range("F24").Value = 344
$ws.Range("F25").Value = 3982
$ws.Range("F27").Value = 42
$ws.Range("F28").Value = 391
$ws.Range("F29").Value = 512
$ws.Range("F30").Value = 1574
$ws.Range("F31").Value = 91
$ws.Range("F33").Value = 2886
$ws.Range("F34").Value = 2063
$ws.Range("F35").Value = 46
$ws.Range("F36").Value = 66
$ws.Range("F37").Value = 106
$ws.Range("F39").Value = 3984
$ws.Range("F40").Value = 391
$ws.Range("F43").Value = 46
$ws.Range("F44").Value = 687
$ws.Range("F45").Value = 1528
